# Apply cryptocurrency price/volume updates as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '63.075.63'
$ws.Cells.Item(2, 5).Value = '  +2.78%  '
$ws.Cells.Item(3, 4).Value = '2.955.83'
$ws.Cells.Item(3, 5).Value = '  +0.96%  '
$ws.Cells.Item(4, 5).Value = '  +0.02%  '
$ws.Cells.Item(5, 4).Value = '''594.74'
$ws.Cells.Item(5, 5).Value = '  -0.46%  '
$ws.Cells.Item(6, 4).Value = '''148.55'
$ws.Cells.Item(6, 5).Value = '  +2.35%  '
$ws.Cells.Item(7, 5).Value = '  -0.06%  '
$ws.Cells.Item(8, 4).Value = '2.954.61'
$ws.Cells.Item(8, 5).Value = '  +1.01%  '
$ws.Cells.Item(9, 5).Value = '  +1.07%  '
$ws.Cells.Item(10, 4).Value = '''7.28'
$ws.Cells.Item(10, 5).Value = '  +4.38%  '
$ws.Cells.Item(11, 5).Value = '  +6.54%  '
$ws.Cells.Item(12, 4).Value = '''0.441'
$ws.Cells.Item(12, 5).Value = '  +0.72%  '
$ws.Cells.Item(13, 5).Value = '  +5.14%  '
$ws.Cells.Item(14, 5).Value = '  -1.69%  '
$ws.Cells.Item(15, 5).Value = '  -0.70%  '
$ws.Cells.Item(16, 4).Value = '3.445.23'
$ws.Cells.Item(17, 4).Value = '63.029.82'
$ws.Cells.Item(17, 5).Value = '  +2.73%  '
$ws.Cells.Item(18, 5).Value = '  +0.31%  '
$ws.Cells.Item(19, 4).Value = '2.939.03'
$ws.Cells.Item(19, 5).Value = '  +0.31%  '
$ws.Cells.Item(20, 4).Value = '''442.60'
$ws.Cells.Item(20, 5).Value = '  +2.53%  '
$ws.Cells.Item(21, 4).Value = '''13.48'
$ws.Cells.Item(21, 5).Value = '  -0.06%  '
$ws.Cells.Item(22, 5).Value = '  -1.03%  '
$ws.Cells.Item(23, 4).Value = '''7.04'
$ws.Cells.Item(23, 5).Value = '  -0.39%  '
$ws.Cells.Item(24, 5).Value = '  +3.53%  '
$ws.Cells.Item(25, 4).Value = '''81.10'
$ws.Cells.Item(25, 5).Value = '  -0.82%  '
$ws.Cells.Item(26, 4).Value = '''2.13'
$ws.Cells.Item(26, 5).Value = '  -2.36%  '
$ws.Cells.Item(27, 4).Value = '''11.77'
$ws.Cells.Item(27, 5).Value = '  +0.42%  '
$ws.Cells.Item(28, 5).Value = '  -0.01%  '
$ws.Cells.Item(29, 5).Value = '  +5.36%  '
$ws.Cells.Item(30, 5).Value = '  -0.44%  '
$ws.Cells.Item(31, 5).Value = '  +0.32%  '
$ws.Cells.Item(33, 2).Value = 'EthereumClassic'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(33, 4).Value = '''26.42'
$ws.Cells.Item(33, 5).Value = '  -0.61%  '
$ws.Cells.Item(34, 2).Value = 'Hedera'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(34, 4).Value = '''0.109'
$ws.Cells.Item(34, 5).Value = '  -0.58%  '
$ws.Cells.Item(35, 5).Value = '  -0.06%  '
$ws.Cells.Item(36, 4).Value = '''0.991'
$ws.Cells.Item(36, 5).Value = '  -1.75%  '
$ws.Cells.Item(37, 4).Value = '''3.15'
$ws.Cells.Item(37, 5).Value = '  +5.72%  '
$ws.Cells.Item(38, 4).Value = '''5.59'
$ws.Cells.Item(38, 5).Value = '  -0.50%  '
$ws.Cells.Item(39, 5).Value = '  +2.89%  '
$ws.Cells.Item(40, 4).Value = '''49.73'
$ws.Cells.Item(40, 5).Value = '  -0.26%  '
$ws.Cells.Item(41, 4).Value = '''8.52'
$ws.Cells.Item(41, 5).Value = '  -0.52%  '
$ws.Cells.Item(42, 5).Value = '  -4.41%  '
$ws.Cells.Item(43, 5).Value = '  +0.29%  '
$ws.Cells.Item(44, 4).Value = '''38.93'
$ws.Cells.Item(44, 5).Value = '  -7.42%  '
$ws.Cells.Item(45, 4).Value = '''135.54'
$ws.Cells.Item(45, 5).Value = '  +1.43%  '
$ws.Cells.Item(46, 4).Value = '2.693.34'
$ws.Cells.Item(46, 5).Value = '  -0.23%  '
$ws.Cells.Item(48, 4).Value = '''360.77'
$ws.Cells.Item(48, 5).Value = '  -0.90%  '
$ws.Cells.Item(50, 5).Value = '  -0.61%  '
$ws.Cells.Item(51, 4).Value = '''22.83'
$ws.Cells.Item(51, 5).Value = '  -3.32%  '
